$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows at position 4 (push existing rows 4-9 down to 6-11) ---
$xlShiftDown = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown
$xlFormatFromLeftOrAbove = [Microsoft.Office.Interop.Excel.XlInsertFormatOrigin]::xlFormatFromLeftOrAbove
$ws.Range("4:5").Insert($xlShiftDown, $xlFormatFromLeftOrAbove)

# --- Row 4: Clostridioides difficile 630 ---
$ws.Range("A4").Value = "Clostridioides_difficile_630"
$ws.Range("B4").Value = "Clostridium difficile 630 (272563.8)"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "Clostridioides difficile 630"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "Clostridioides difficile"
$ws.Range("F4").Value = "Clostridioides"
$ws.Range("G4").Value = "Peptostreptococcaceae"
$ws.Range("H4").Value = "Clostridiales"
$ws.Range("I4").Value = "Clostridia"
$ws.Range("J4").Value = "Firmicutes"
$ws.Range("K4").Value = 272563
$ws.Range("L4").Value = 272563.8
$ws.Range("M4").Value = "Gram+"
$ws.Range("N4").Value = "Yes"

# --- Row 5: Lactobacillus jensenii 269-3 ---
$ws.Range("A5").Value = "Lactobacillus_jensenii_269_3"
$ws.Range("B5").Value = "Lactobacillus jensenii 269-3 (596325.3)"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "Lactobacillus jensenii 269-3"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "Lactobacillus jensenii"
$ws.Range("F5").Value = "Lactobacillus"
$ws.Range("G5").Value = "Lactobacillaceae"
$ws.Range("H5").Value = "Lactobacillales"
$ws.Range("I5").Value = "Bacilli"
$ws.Range("J5").Value = "Firmicutes"
$ws.Range("K5").Value = 596325
$ws.Range("L5").Value = 596325.3
$ws.Range("M5").Value = "Gram+"
$ws.Range("N5").Value = "Yes"

# --- Re-apply the sort state so it now refers to the shifted block (A6:J10) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A6:A10")) | Out-Null
$ws.Sort.SetRange($ws.Range("A6:J10"))
$ws.Sort.Apply()

# --- Split the old merged "C:D" column-width definition into two columns ---
# (column D keeps its original width; column C gets a new, narrower width)
$ws.Columns.Item(3).ColumnWidth = 6.5
$ws.Columns.Item(4).ColumnWidth = 30.5

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("B18").Select()

Write-Output "edit applied"
